$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 90480

# Row 3
$ws.Range("A3").Value = 112227891
$ws.Range("B3").Value = 89549
$ws.Range("E3").Value = 1108
$ws.Range("F3").Value = "Harticka"
$ws.Range("G3").Value = "Pelloporus leporinus"
$ws.Range("H3").Value = "(Fr.) Krieglst."
$ws.Range("Q3").Value = 496245
$ws.Range("R3").Value = 6934459
$ws.Range("Z3").Value = "11:32"
$ws.Range("AB3").Value = "11:32"

# Row 4
$ws.Range("A4").Value = 112228055
$ws.Range("B4").Value = 89517
$ws.Range("D4").Value = "LC"
$ws.Range("E4").Value = 5447
$ws.Range("F4").Value = "Vedticka"
$ws.Range("G4").Value = "Fuscoporia viticola"
$ws.Range("H4").Value = "(Schwein.) Murrill"
$ws.Range("Q4").Value = 496258
$ws.Range("R4").Value = 6934460
$ws.Range("Z4").Value = "11:46"
$ws.Range("AB4").Value = "11:46"

# Row 5
$ws.Range("A5").Value = 112228190
$ws.Range("B5").Value = 89571
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 5432
$ws.Range("F5").Value = "Granticka"
$ws.Range("G5").Value = "Porodaedalea chrysoloma"
$ws.Range("H5").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q5").Value = 496305
$ws.Range("R5").Value = 6934462
$ws.Range("Z5").Value = "11:47"
$ws.Range("AB5").Value = "11:47"

# Row 6
$ws.Range("A6").Value = 112228201
$ws.Range("B6").Value = 89553
$ws.Range("E6").Value = 1202
$ws.Range("F6").Value = "Ullticka"
$ws.Range("G6").Value = "Phellinidium ferrugineofuscum"
$ws.Range("H6").Value = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("Q6").Value = 496302
$ws.Range("R6").Value = 6934437
$ws.Range("Z6").Value = "11:57"
$ws.Range("AB6").Value = "11:57"

# Row 7
$ws.Range("A7").Value = 112227657
$ws.Range("B7").Value = 81385
$ws.Range("E7").Value = 1312
$ws.Range("F7").Value = "Gammelgransskål"
$ws.Range("G7").Value = "Pseudographis pinicola"
$ws.Range("H7").Value = "(Nyl.) Rehm"
$ws.Range("Q7").Value = 496238
$ws.Range("R7").Value = 6934504
$ws.Range("Z7").Value = "11:20"
$ws.Range("AB7").Value = "11:20"

# Row 8
$ws.Range("B8").Value = 90213
